$wb = $excel.ActiveWorkbook

# --- Update AddAdminData test data (row 3): replace the "James Butler" entry
# with a new "Mohai Islam" employee / username / password entry ---
$wsAdmin = $wb.Worksheets.Item("AddAdminData")
$wsAdmin.Cells.Item(3, 1).Value = "Mohai Islam"
$wsAdmin.Cells.Item(3, 2).Value = "mohai1"
$wsAdmin.Cells.Item(3, 3).Value = "mohai123"
$wsAdmin.Cells.Item(3, 4).Value = "mohai123"

# --- Leave the workbook with AddAdminData as the active sheet / tab, with
# cell B3 selected (matches the saved state captured in the workbook) ---
$wsAdmin.Activate()
$wsAdmin.Range("B3").Select()
